$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("M")

# Replace specific yearly figures (column J = oldest year, 2012) with "NA"
# as the data was unavailable / not reported for that period.
$ws.Range("J21").Value = "NA"   # Earnings Before Interest And Taxes
$ws.Range("J83").Value = "NA"   # Depreciation
$ws.Range("J94").Value = "NA"   # Total Cash Flows From Investing Activities
$ws.Range("J100").Value = "NA"  # Total Cash Flows From Financing Activities
$ws.Range("J101").Value = "NA"  # Effect Of Exchange Rate Changes

# Update "Capital Expenditures" row (row 91) with revised figures
$ws.Range("D91").Value = -487000
$ws.Range("E91").Value = -596000
$ws.Range("F91").Value = -777000
$ws.Range("G91").Value = -770000
$ws.Range("H91").Value = -607000
$ws.Range("I91").Value = -698000
$ws.Range("J91").Value = -555000

$wb.Save()
